# Updates cryptos list (prices / 1h volume %) and swaps the dogwifhat /
# Fetch.AI rows (45-46), matching the "Updated cryptos list" GitHub Actions
# commit. Values are prefixed with a literal leading apostrophe so Excel
# keeps them as text (matching the original inline-string cells, e.g.
# "67.052.96" or "0.0430") instead of silently coercing them to numbers;
# the Style is then reset to "Normal" so no stray text-format style is
# left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''67.052.96'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '''  +3.45%  '
$ws.Range('E2').Style = "Normal"
$ws.Range('D3').Value = '''3.452.47'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '''  +2.95%  '
$ws.Range('E3').Style = "Normal"
$ws.Range('E4').Value = '''  -0.03%  '
$ws.Range('E4').Style = "Normal"
$ws.Range('D5').Value = '''580.56'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '''  +4.76%  '
$ws.Range('E5').Style = "Normal"
$ws.Range('D6').Value = '''186.47'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '''  +7.27%  '
$ws.Range('E6').Style = "Normal"
$ws.Range('D7').Value = '''0.632'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '''  +0.47%  '
$ws.Range('E7').Style = "Normal"
$ws.Range('D8').Value = '''3.444.56'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '''  +2.97%  '
$ws.Range('E8').Style = "Normal"
$ws.Range('E9').Value = '''  +0.01%  '
$ws.Range('E9').Style = "Normal"
$ws.Range('E10').Value = '''  -0.93%  '
$ws.Range('E10').Style = "Normal"
$ws.Range('D11').Value = '''0.646'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '''  +1.51%  '
$ws.Range('E11').Style = "Normal"
$ws.Range('D12').Value = '''56.12'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '''  +4.69%  '
$ws.Range('E12').Style = "Normal"
$ws.Range('E13').Value = '''  -1.26%  '
$ws.Range('E13').Style = "Normal"
$ws.Range('D14').Value = '''9.39'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '''  +3.18%  '
$ws.Range('E14').Style = "Normal"
$ws.Range('D15').Value = '''3.999.13'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '''  +2.94%  '
$ws.Range('E15').Style = "Normal"
$ws.Range('D16').Value = '''18.70'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '''  +2.81%  '
$ws.Range('E16').Style = "Normal"
$ws.Range('D17').Value = '''3.447.80'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '''  +3.09%  '
$ws.Range('E17').Style = "Normal"
$ws.Range('D18').Value = '''67.008.98'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '''  +3.68%  '
$ws.Range('E18').Style = "Normal"
$ws.Range('D19').Value = '''12.10'
$ws.Range('D19').Style = "Normal"
$ws.Range('E20').Value = '''  -2.31%  '
$ws.Range('E20').Style = "Normal"
$ws.Range('E21').Value = '''  +3.35%  '
$ws.Range('E21').Style = "Normal"
$ws.Range('D22').Value = '''487.19'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '''  +8.64%  '
$ws.Range('E22').Style = "Normal"
$ws.Range('D23').Value = '''5.31'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '''  +7.95%  '
$ws.Range('E23').Style = "Normal"
$ws.Range('D24').Value = '''16.91'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '''  +24.24%  '
$ws.Range('E24').Style = "Normal"
$ws.Range('D25').Value = '''4.37'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '''  +7.86%  '
$ws.Range('E25').Style = "Normal"
$ws.Range('D26').Value = '''89.61'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '''  +3.40%  '
$ws.Range('E26').Style = "Normal"
$ws.Range('E27').Value = '''  +3.17%  '
$ws.Range('E27').Style = "Normal"
$ws.Range('D28').Value = '''10.93'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '''  +2.20%  '
$ws.Range('E28').Style = "Normal"
$ws.Range('D29').Value = '''9.05'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '''  +4.93%  '
$ws.Range('E29').Style = "Normal"
$ws.Range('D30').Value = '''31.39'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '''  +1.59%  '
$ws.Range('E30').Style = "Normal"
$ws.Range('D31').Value = '''7.23'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '''  +10.64%  '
$ws.Range('E31').Style = "Normal"
$ws.Range('D32').Value = '''598.05'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '''  +3.80%  '
$ws.Range('E32').Style = "Normal"
$ws.Range('D33').Value = '''11.76'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '''  +2.95%  '
$ws.Range('E33').Style = "Normal"
$ws.Range('D34').Value = '''63.86'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '''  +1.20%  '
$ws.Range('E34').Style = "Normal"
$ws.Range('D35').Value = '''0.112'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '''  +4.23%  '
$ws.Range('E35').Style = "Normal"
$ws.Range('E36').Value = '''  +6.91%  '
$ws.Range('E36').Style = "Normal"
$ws.Range('E37').Value = '''  -0.04%  '
$ws.Range('E37').Style = "Normal"
$ws.Range('D38').Value = '''36.65'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '''  +2.95%  '
$ws.Range('E38').Style = "Normal"
$ws.Range('E39').Value = '''  +3.95%  '
$ws.Range('E39').Style = "Normal"
$ws.Range('D40').Value = '''3.53'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '''  -3.28%  '
$ws.Range('E40').Style = "Normal"
$ws.Range('D41').Value = '''3.256.43'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '''  +5.97%  '
$ws.Range('E41').Style = "Normal"
$ws.Range('E42').Value = '''  +1.93%  '
$ws.Range('E42').Style = "Normal"
$ws.Range('D43').Value = '''2.91'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '''  +6.36%  '
$ws.Range('E43').Style = "Normal"
$ws.Range('D44').Value = '''0.0430'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '''  +3.30%  '
$ws.Range('E44').Style = "Normal"
$ws.Range('B45').Value = '''Fetch.AI'
$ws.Range('B45').Style = "Normal"
$ws.Range('C45').Value = '''https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('C45').Style = "Normal"
$ws.Range('D45').Value = '''2.52'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '''  +3.04%  '
$ws.Range('E45').Style = "Normal"
$ws.Range('B46').Value = '''dogwifhat'
$ws.Range('B46').Style = "Normal"
$ws.Range('C46').Value = '''https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('C46').Style = "Normal"
$ws.Range('D46').Value = '''2.77'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '''  +22.84%  '
$ws.Range('E46').Style = "Normal"
$ws.Range('D47').Value = '''3.22'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '''  +1.49%  '
$ws.Range('E47').Style = "Normal"
$ws.Range('E48').Value = '''  +0.31%  '
$ws.Range('E48').Style = "Normal"
$ws.Range('E49').Value = '''  +13.54%  '
$ws.Range('E49').Style = "Normal"
$ws.Range('E50').Value = '''  +5.59%  '
$ws.Range('E50').Style = "Normal"
$ws.Range('E51').Value = '''  +0.06%  '
$ws.Range('E51').Style = "Normal"
